# thêm điểm iteration 3
# Create the new "Iteration 3" worksheet by duplicating the "Iteration 2" sheet
# (same layout/template) and then updating the group scores / comments to the
# values for Iteration 3.

$wb = $excel.ActiveWorkbook

$srcSheet = $wb.Worksheets.Item("Iteration 2, Specs and Design")

# Duplicate the template sheet and place the copy right after it (i.e. at the
# end of the workbook), mirroring how the new "Iteration 3" tab was added.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$srcSheet.Copy($null, $lastSheet)

$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "Iteration 3"

# Header title for the new iteration
$ws.Range("B1").Value = "Iteration 3"

# ---- Nhóm 21 (rows 4-9) : unchanged ----

# ---- Nhóm 22 (rows 10-15) ----
$ws.Range("E10").Value = 40
$ws.Range("F10").Value = "lần insert thứ nhất các bạn làm csdl sai mà không xem xét trước khi nộp "
$ws.Range("D11").Value = 8
$ws.Range("D12").Value = 8
$ws.Range("D13").Value = 8
$ws.Range("D14").Value = 8
$ws.Range("D15").Value = 8

# ---- Nhóm 23 (rows 16-21) ----
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = "nộp đúng giờ và đầy đủ"
$ws.Range("D17").Value = 10
$ws.Range("D18").Value = 10
$ws.Range("D19").Value = 10
$ws.Range("D20").Value = 10
$ws.Range("D21").Value = 10

# ---- Nhóm 24 (rows 22-27) ----
$ws.Range("E22").Value = 50
$ws.Range("F22").Value = "nộp đúng giờ và đầy đủ"
$ws.Range("D23").Value = 10
$ws.Range("D24").Value = 10
$ws.Range("D25").Value = 10
$ws.Range("D26").Value = 10
$ws.Range("D27").Value = 10

# ---- Nhóm 25 (rows 28-33) : unchanged ----

# ---- Nhóm 26 (rows 34-39) ----
$ws.Range("E34").Value = 50
$ws.Range("F34").Value = "các bạn còn làm một số chức năng có sai sót  như mình đã đưa ra, nhưng nhìn chung các bạn có ý thức sửa và làm lại tốt"
$ws.Range("D35").Value = 10
$ws.Range("D36").Value = 10
$ws.Range("D37").Value = 10
$ws.Range("D38").Value = 10
$ws.Range("D39").Value = 10

# ---- Nhóm 27 (rows 40-44) ----
$ws.Range("E40").Value = 30
$ws.Range("F40").Value = "các bạn insert lần 2 add sai dữ liệu"
$ws.Range("D41").Value = 8
$ws.Range("D42").Value = 8
$ws.Range("D43").Value = 7
$ws.Range("D44").Value = 7

# ---- Nhóm 28 (rows 45-50) : unchanged ----

# ---- Nhóm 29 (rows 51-56) ----
$ws.Range("E51").Value = 50
$ws.Range("F51").Value = "nộp đúng giờ và đầy đủ"
$ws.Range("D52").Value = 10
$ws.Range("D53").Value = 10
$ws.Range("D54").Value = 10
$ws.Range("D55").Value = 10
$ws.Range("D56").Value = 10

# ---- Nhóm 44 (row 57-62) ----
$ws.Range("E57").Value = 0
$ws.Range("F57").Value = "các bạn liên tục không nộp bài hoặc không tổng hợp dữ liệu."
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 0
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 0
$ws.Range("D62").Value = 0

# Leave the previous (Iteration 2) sheet's selection resting on the full
# used range, then make the new "Iteration 3" sheet active with F16 selected,
# matching the final view state of the edit.
$srcSheet.Select()
$srcSheet.Range("A1:F62").Select()

$ws.Select()
$ws.Range("F16").Select()
